$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4736.5625
$ws.Range("J116").Value = 8657.143
$ws.Range("L116").Value = 8657.143
$ws.Range("N116").Value = -15541.143
$ws.Range("H129").Value = 848.36536
$ws.Range("I129").Value = 471.53845
$ws.Range("J129").Value = 973.97437
$ws.Range("K129").Value = 1414.61535
$ws.Range("L129").Value = 2921.92311
$ws.Range("M129").Value = 3585.38465
$ws.Range("N129").Value = -12921.92311
$ws.Range("H131").Value = 775.6316
$ws.Range("I131").Value = 587.7
$ws.Range("J131").Value = 984.44446
$ws.Range("K131").Value = 1763.1
$ws.Range("L131").Value = 2953.33338
$ws.Range("M131").Value = 3276.9
$ws.Range("N131").Value = -13033.33338
$ws.Range("H135").Value = 40306.348
$ws.Range("I135").Value = 72987.36
$ws.Range("J135").Value = 2178.5
$ws.Range("K135").Value = 656886.24
$ws.Range("L135").Value = 19606.5
$ws.Range("M135").Value = -654351.24
$ws.Range("N135").Value = -24676.5
$ws.Range("H138").Value = 2018.4318
$ws.Range("I138").Value = 1433.3572
$ws.Range("J138").Value = 2129.1216
$ws.Range("K138").Value = 4300.071599999999
$ws.Range("L138").Value = 6387.364799999999
$ws.Range("M138").Value = 839.9284000000007
$ws.Range("N138").Value = -16667.3648

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1729.2632
$ws.Range("I2").Value = 1649.4615
$ws.Range("J2").Value = 1902.1666
$ws.Range("K2").Value = 1649.4615
$ws.Range("L2").Value = 1902.1666
$ws.Range("M2").Value = -1536.4615
$ws.Range("N2").Value = -2128.1666
$ws.Range("H32").Value = 7601988.5
$ws.Range("I32").Value = 8801018
$ws.Range("J32").Value = 8135.5
$ws.Range("K32").Value = 8801018
$ws.Range("L32").Value = 8135.5
$ws.Range("M32").Value = -8800731
$ws.Range("N32").Value = -8709.5
$ws.Range("H64").Value = 21792.285
$ws.Range("J64").Value = 21930.154
$ws.Range("L64").Value = 21930.154
$ws.Range("N64").Value = -22426.154
$ws.Range("H67").Value = 21792.285
$ws.Range("J67").Value = 21930.154
$ws.Range("L67").Value = 21930.154
$ws.Range("N67").Value = -23646.154
$ws.Range("H74").Value = 14001323
$ws.Range("I74").Value = 20917772
$ws.Range("K74").Value = 20917772
$ws.Range("M74").Value = -20916898
$ws.Range("H77").Value = 14001323
$ws.Range("I77").Value = 20917772
$ws.Range("K77").Value = 104588860
$ws.Range("M77").Value = -104584492
$ws.Range("H116").Value = 1729.2632
$ws.Range("I116").Value = 1649.4615
$ws.Range("J116").Value = 1902.1666
$ws.Range("K116").Value = 1649.4615
$ws.Range("L116").Value = 1902.1666
$ws.Range("M116").Value = 644.5385000000001
$ws.Range("N116").Value = -6490.1666
$ws.Range("H132").Value = 66351.35000000001
$ws.Range("I132").Value = 39960.77
$ws.Range("K132").Value = 119882.31
$ws.Range("M132").Value = -117352.31

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1729.2632
$ws.Range("I3").Value = 1649.4615
$ws.Range("J3").Value = 1902.1666
$ws.Range("K3").Value = 1649.4615
$ws.Range("L3").Value = 1902.1666
$ws.Range("M3").Value = -1535.4615
$ws.Range("N3").Value = -2130.1666
$ws.Range("H62").Value = 40181
$ws.Range("J62").Value = 40181
$ws.Range("L62").Value = 40181
$ws.Range("N62").Value = -41553
$ws.Range("H65").Value = 40181
$ws.Range("J65").Value = 40181
$ws.Range("L65").Value = 120543
$ws.Range("N65").Value = -127407
$ws.Range("H94").Value = 950.35297
$ws.Range("I94").Value = 797.34485
$ws.Range("J94").Value = 1837.8
$ws.Range("K94").Value = 797.34485
$ws.Range("L94").Value = 1837.8
$ws.Range("M94").Value = -346.34485
$ws.Range("N94").Value = -2739.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2246.4614
$ws.Range("I16").Value = 2034.7778
$ws.Range("K16").Value = 2034.7778
$ws.Range("M16").Value = -1747.7778
$ws.Range("H31").Value = 5000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 5000
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -5590
$ws.Range("H34").Value = 5000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5000
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -5404
$ws.Range("H113").Value = 2246.4614
$ws.Range("I113").Value = 2034.7778
$ws.Range("K113").Value = 2034.7778
$ws.Range("M113").Value = 135.2221999999999
$ws.Range("H122").Value = 1294.5333
$ws.Range("I122").Value = 1055.4286
$ws.Range("J122").Value = 2131.4
$ws.Range("K122").Value = 3166.2858
$ws.Range("L122").Value = 6394.200000000001
$ws.Range("M122").Value = -716.2857999999997
$ws.Range("N122").Value = -11294.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 639.8889
$ws.Range("I92").Value = 260
$ws.Range("J92").Value = 829.8333
$ws.Range("K92").Value = 780
$ws.Range("L92").Value = 2489.4999
$ws.Range("M92").Value = 468
$ws.Range("N92").Value = -4985.4999
$ws.Range("H122").Value = 1131.0625
$ws.Range("I122").Value = 448.75
$ws.Range("J122").Value = 1358.5
$ws.Range("K122").Value = 4038.75
$ws.Range("L122").Value = 12226.5
$ws.Range("M122").Value = -1588.75
$ws.Range("N122").Value = -17126.5
$ws.Range("H132").Value = 1472.2778
$ws.Range("I132").Value = 1283.6666
$ws.Range("J132").Value = 1660.8889
$ws.Range("K132").Value = 11552.9994
$ws.Range("L132").Value = 14948.0001
$ws.Range("M132").Value = -9022.999400000001
$ws.Range("N132").Value = -20008.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2828.7273
$ws.Range("I122").Value = 2012.5
$ws.Range("J122").Value = 5005.3335
$ws.Range("K122").Value = 6037.5
$ws.Range("L122").Value = 15016.0005
$ws.Range("M122").Value = -3587.5
$ws.Range("N122").Value = -19916.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2639.8667
$ws.Range("I7").Value = 2549.9167
$ws.Range("J7").Value = 2999.6667
$ws.Range("K7").Value = 2549.9167
$ws.Range("L7").Value = 2999.6667
$ws.Range("M7").Value = -2437.9167
$ws.Range("N7").Value = -3223.6667
$ws.Range("H126").Value = 2639.8667
$ws.Range("I126").Value = 2549.9167
$ws.Range("J126").Value = 2999.6667
$ws.Range("K126").Value = 7649.750100000001
$ws.Range("L126").Value = 8999.000100000001
$ws.Range("M126").Value = -5179.750100000001
$ws.Range("N126").Value = -13939.0001
$ws.Range("H136").Value = 252698.88
$ws.Range("I136").Value = 252550.25
$ws.Range("K136").Value = 757650.75
$ws.Range("M136").Value = -755100.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1120.6875
$ws.Range("I113").Value = 696.6667
$ws.Range("J113").Value = 1375.1
$ws.Range("K113").Value = 2090.0001
$ws.Range("L113").Value = 4125.299999999999
$ws.Range("M113").Value = 79.9998999999998
$ws.Range("N113").Value = -8465.299999999999
$ws.Range("H122").Value = 1837.5518
$ws.Range("I122").Value = 1521
$ws.Range("J122").Value = 3357
$ws.Range("K122").Value = 4563
$ws.Range("L122").Value = 10071
$ws.Range("M122").Value = -2113
$ws.Range("N122").Value = -14971
$ws.Range("H136").Value = 128968.375
$ws.Range("I136").Value = 95363.09
$ws.Range("J136").Value = 202900
$ws.Range("K136").Value = 286089.27
$ws.Range("L136").Value = 608700
$ws.Range("M136").Value = -283539.27
$ws.Range("N136").Value = -613800
